# "Generate Report for Handoff"
#
# The localization-status workbook tracks per-file handoff status across
# three sheets: "Overview", "zh-cn" and "de-de". A new handoff run moved the
# file "b8e29229-0a76-4ad9-8eff-ecbab82f272d.md" (previously "In Translation",
# sitting in row 6) out to translators -- it is now "Ready for handoff" with a
# fresh handoff datetime, and "e76e6db7-09bf-43d7-a250-ee755d8d0490.md"
# (previously row 7) takes over row 6 still "In Translation".
# In effect rows 6 and 7 swap identity/content on every sheet, and the file
# that moves to row 7 picks up its new "Ready for handoff" status + handoff
# timestamp.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A (File Name), B (zh-cn), C (de-de), D (Latest
# Handoff Date)
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A6").Value = "e76e6db7-09bf-43d7-a250-ee755d8d0490.md"
$ov.Range("B6").Value = "In Translation"
$ov.Range("C6").Value = "In Translation"
$ov.Range("D6").Value = "2016-32-12 20:32:35"

$ov.Range("A7").Value = "b8e29229-0a76-4ad9-8eff-ecbab82f272d.md"
$ov.Range("B7").Value = "Ready for handoff"
$ov.Range("C7").Value = "Ready for handoff"
$ov.Range("D7").Value = "2016-35-12 20:35:42"

# ---------------------------------------------------------------------
# Sheet "zh-cn": columns A (Source File Name), B (File Extension),
# C (Status), D (Latest Handoff File), E (Latest Handoff Datetime)
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A6").Value = "e76e6db7-09bf-43d7-a250-ee755d8d0490.md"
$zh.Range("C6").Value = "In Translation"
$zh.Range("D6").Value = "e76e6db7-09bf-43d7-a250-ee755d8d0490.9c209513f88ec7d62cc6042e3ff61e636822feaf.zh-cn.xlf"
$zh.Range("E6").Value = "2016-03-12 20:32:32"

$zh.Range("A7").Value = "b8e29229-0a76-4ad9-8eff-ecbab82f272d.md"
$zh.Range("C7").Value = "Ready for handoff"
$zh.Range("D7").Value = "b8e29229-0a76-4ad9-8eff-ecbab82f272d.b7d9c615aab95b7f1f833faf2ea1bbad81e48d6f.zh-cn.xlf"
$zh.Range("E7").Value = "2016-03-12 20:35:39"

# ---------------------------------------------------------------------
# Sheet "de-de": columns A (Source File Name), B (File Extension),
# C (Status), D (Latest Handoff File), E (Latest Handoff Datetime)
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A6").Value = "e76e6db7-09bf-43d7-a250-ee755d8d0490.md"
$de.Range("C6").Value = "In Translation"
$de.Range("D6").Value = "e76e6db7-09bf-43d7-a250-ee755d8d0490.9c209513f88ec7d62cc6042e3ff61e636822feaf.de-de.xlf"
$de.Range("E6").Value = "2016-03-12 20:32:35"

$de.Range("A7").Value = "b8e29229-0a76-4ad9-8eff-ecbab82f272d.md"
$de.Range("C7").Value = "Ready for handoff"
$de.Range("D7").Value = "b8e29229-0a76-4ad9-8eff-ecbab82f272d.b7d9c615aab95b7f1f833faf2ea1bbad81e48d6f.de-de.xlf"
$de.Range("E7").Value = "2016-03-12 20:35:42"
